$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.772.97"
$ws.Range("E2").Value = "  +2.22%  "

$ws.Range("D3").Value = "3.036.84"
$ws.Range("E3").Value = "  +1.52%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'512.87"
$ws.Range("E5").Value = "  +1.26%  "

$ws.Range("D6").Value = "'140.73"
$ws.Range("E6").Value = "  +3.19%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +2.50%  "

$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").Value = "'0.112"
$ws.Range("E10").Value = "  +2.84%  "

$ws.Range("D11").Value = "'0.369"
$ws.Range("E11").Value = "  +4.60%  "

$ws.Range("D12").Value = "3.552.10"
$ws.Range("E12").Value = "  +1.43%  "

$ws.Range("E13").Value = "  +1.49%  "

$ws.Range("D14").Value = "'26.77"
$ws.Range("E14").Value = "  +3.78%  "

$ws.Range("D15").Value = "'0.0000167"
$ws.Range("E15").Value = "  +8.52%  "

$ws.Range("D16").Value = "57.773.09"
$ws.Range("E16").Value = "  +2.20%  "

$ws.Range("D17").Value = "'6.26"
$ws.Range("E17").Value = "  +7.79%  "

$ws.Range("D18").Value = "3.034.62"
$ws.Range("E18").Value = "  +1.63%  "

$ws.Range("D19").Value = "'12.90"
$ws.Range("E19").Value = "  +3.34%  "

$ws.Range("D20").Value = "'8.05"
$ws.Range("E20").Value = "  +3.11%  "

$ws.Range("D21").Value = "'333.63"
$ws.Range("E21").Value = "  +2.13%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").Value = "'0.501"
$ws.Range("E23").Value = "  +5.50%  "

$ws.Range("D24").Value = "'64.80"
$ws.Range("E24").Value = "  +4.01%  "

$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").Value = "0.0₃0939"
$ws.Range("E27").Value = "  +3.85%  "

$ws.Range("D28").Value = "'6.81"
$ws.Range("E28").Value = "  +5.18%  "

$ws.Range("E29").Value = "  +8.31%  "

$ws.Range("E30").Value = "  +2.86%  "

$ws.Range("D31").Value = "'1.22"
$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("D32").Value = "'20.82"
$ws.Range("E32").Value = "  +1.09%  "

$ws.Range("D33").Value = "'4.75"
$ws.Range("E33").Value = "  +5.87%  "

$ws.Range("D34").Value = "'155.53"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").Value = "'5.91"
$ws.Range("E35").Value = "  +5.83%  "

$ws.Range("D36").Value = "'1.29"
$ws.Range("E36").Value = "  +1.88%  "

$ws.Range("D37").Value = "'24.93"
$ws.Range("E37").Value = "  +6.20%  "

$ws.Range("E38").Value = "  +1.65%  "

$ws.Range("D39").Value = "3.069.03"
$ws.Range("E39").Value = "  +1.45%  "

$ws.Range("D40").Value = "'37.53"
$ws.Range("E40").Value = "  +3.27%  "

$ws.Range("D41").Value = "'3.90"
$ws.Range("E41").Value = "  +8.44%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").Value = "2.319.30"
$ws.Range("E43").Value = "  +2.63%  "

$ws.Range("D44").Value = "'0.658"
$ws.Range("E44").Value = "  +2.22%  "

$ws.Range("E45").Value = "  +1.67%  "

$ws.Range("D46").Value = "'1.00"
$ws.Range("E46").Value = "  +1.06%  "

$ws.Range("D47").Value = "'6.06"
$ws.Range("E47").Value = "  +5.02%  "

$ws.Range("D48").Value = "'0.0242"
$ws.Range("E48").Value = "  +2.28%  "

$ws.Range("D49").Value = "'19.68"
$ws.Range("E49").Value = "  +2.75%  "

$ws.Range("D50").Value = "'1.87"
$ws.Range("E50").Value = "  -4.79%  "

$ws.Range("D51").Value = "'0.0896"
$ws.Range("E51").Value = "  +3.08%  "
